# Change the table style on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") from
# the built-in style {0885E01A-58E2-4A35-B4CB-528AEB7FC80C} to
# {599ECA78-6ABE-4091-B448-65C46F360636}.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

foreach ($sh in $s.Shapes) {
    if ($sh.HasTable) {
        $tbl = $sh.Table
        $tbl.ApplyStyle("{599ECA78-6ABE-4091-B448-65C46F360636}")
    }
}
